# Sprint 1 Meeting 8 (mistake in last commit, that was meeting 7)
# Fill in the "3/10/2019: 10pm" (column J) answers for all four team members.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "Got intents working and linked all of the screens I have built thus far. Also renamed layout components in the format of Activity_NameType. This improved the readability in the Java code when mapping java objets to their corrects IDs."
$ws.Range("J3").Value = "I worked on the XML and Java file for the main menu, sign up page, and about page."
$ws.Range("J4").Value = "Created the layout and some functionality for the review page"
$ws.Range("J5").Value = "Continued attempting to port a Unity project to Android studio"

$ws.Range("J6").Value = "Finish planning all topics of discussion for our 'merge meeting', update the github, and continue to look into firebase"
$ws.Range("J7").Value = "I will work on my assigned issues. Every Issue assigned to me as well as the design and layout of every button and image in all pages for all issues of all assignee."
$ws.Range("J8").Value = "Finish functionality of review, help, and info pages"
$ws.Range("J9").Value = "Continue looking into porting a Unity project to Android Studio"

$ws.Range("J10").Value = "Not currently"
$ws.Range("J11").Value = "No, I am still working with all cylinders pumping"
$ws.Range("J12").Value = "Travel "
$ws.Range("J13").Value = "Several errors have prevented me from successfully running a Unity proect in Android Studio"

$ws.Range("J14").Value = "The use of intentions"
$ws.Range("J15").Value = "Working together is better than alone!"
$ws.Range("J16").Value = "Learning different ways of adjusting EditText and Ratings"
$ws.Range("J17").Value = "Learned more about how to run a Unity project in Andriod Studio"

$ws.Range("J18").Value = "Changing the format of layout widget names will need to happen at some point to standardize them"
$ws.Range("J19").Value = "A consistent color scheme/theme will be a necessary addition. However, most likely will not occur until the end of the project"
$ws.Range("J20").Value = "Not yet"
$ws.Range("J21").Value = "No changes currently have to be made to the current plan for the project"

# Row heights grew because the new text in column J wraps onto more lines.
$ws.Rows.Item(7).RowHeight = 97
$ws.Rows.Item(19).RowHeight = 82

# Move the view / selection to where the new answers were entered.
$ws.Range("J19").Select()

# Reflect the scrolled viewport / resized window from the original edit session.
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 8
$win.Left = -80
$win.Top = 0
$win.Width = 10290
$win.Height = 7360
